# Redux.docx edit: split the closing paragraph (the one ending in
# " diyoruz." and carrying the trailing _GoBack bookmark) so that two
# new paragraphs are inserted between the existing text and the
# bookmark:
#   1. "Özet : ..." summary paragraph
#   2. "Componentin Reducera Bağlanması" bold heading paragraph
# The _GoBack bookmark ends up alone in its own (now last) paragraph,
# exactly as in the target document.

$d = $word.ActiveDocument

# Find the end of the " diyoruz." run; that is where the split must
# happen (immediately before the _GoBack bookmark that currently sits
# at the very end of the last paragraph).
$rng = $d.Content
$found = $rng.Find.Execute(" diyoruz.", $false, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$insertPos = $rng.Start

# Insert three paragraph breaks right before the bookmark. The first
# break pushes the bookmark into its own trailing paragraph; the next
# two create a pair of empty paragraphs ahead of it, ready to receive
# the new content.
$d.Range($insertPos, $insertPos).InsertBefore("`r")
$d.Range($insertPos, $insertPos).InsertBefore("`r")
$d.Range($insertPos, $insertPos).InsertBefore("`r")

$paraCount = $d.Paragraphs.Count

# Second-to-last paragraph -> summary text.
$summaryPara = $d.Paragraphs($paraCount - 2)
$summaryPara.Range.Text = "Özet : İlk başta npm ile redux ve react-redux kur.Burada redux ,redux ın kendisidir.React-redux ise componentlerimizi redux a bağlamak için kullanılan küyüphanedir."

# Paragraph right before the bookmark paragraph -> bold heading.
$headingPara = $d.Paragraphs($paraCount - 1)
$headingPara.Range.Text = "Componentin Reducera Bağlanması"
$headingPara.Range.Bold = 1
